# Slide 8 ("Comprehensive Literature Search"), content placeholder shape:
# paragraph 4 ("41 Trials: citation-exports(CochraneTrials) (2).ris") gets
# two small text tweaks:
#   "citation-exports("  -> "citation-exports ("   (space added before paren)
#   ") (2)."             -> ")."                    (" (2)" removed)
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the "41 Trials: ..." paragraph defensively (rather than assuming
# it is always paragraph 4) by scanning for the distinctive text.
$paraIndex = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    if ($tr.Paragraphs($i).Text -like "*citation-exports(*") {
        $paraIndex = $i
        break
    }
}

$para   = $tr.Paragraphs($paraIndex)
$pStart = $para.Start
$pText  = $para.Text

# 1) "citation-exports(" -> "citation-exports ("
$oldD = "citation-exports("
$posD = $pText.IndexOf($oldD)
$runD = $tr.Characters($pStart + $posD, $oldD.Length)
$runD.Text = "citation-exports ("

# Re-fetch paragraph/text after the first mutation shifted later offsets.
$para   = $tr.Paragraphs($paraIndex)
$pStart = $para.Start
$pText  = $para.Text

# 2) ") (2)." -> ")."
$oldF = ") (2)."
$posF = $pText.IndexOf($oldF)
$runF = $tr.Characters($pStart + $posF, $oldF.Length)
$runF.Text = ")."
